# Stats.xlsx / StatDefinition sheet update
# - Inserts 3 new columns (isPrimaryStat, isSecondaryStat, isTertiaryStat) before the
#   old "valueMultiplier" column (old Q -> new T, old R -> new U, old S -> new V)
# - Flags each existing stat row as primary/secondary/tertiary via the new columns
# - Inserts a new "statMult" definition row before the old "gold" row
# - Appends a new "rarity" definition row at the end

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert the three new columns before the old column Q (valueMultiplier) ---
$ws.Range("Q1:S1").EntireColumn.Insert()

# --- 2. Insert a new row before the old row 38 ("gold") ---
$ws.Range("A38:A38").EntireRow.Insert()

# --- 3. New header cells for the inserted columns ---
$ws.Range("Q1").Value = "isPrimaryStat"
$ws.Range("R1").Value = "isSecondaryStat"
$ws.Range("S1").Value = "isTertiaryStat"

# --- 4. Mark each stat as primary / secondary / tertiary ---

# Secondary stats (column R)
$secondaryRows = @(3,4,5,6,8,9,10,11,12,13,14,16,18,20,21,22,23,24,25,26,30,31,32,33,34,35)
foreach ($r in $secondaryRows) {
    $ws.Range("R$r").Value = $true
}

# Primary stats (column Q)
$primaryRows = @(15,17,19)
foreach ($r in $primaryRows) {
    $ws.Range("Q$r").Value = $true
}

# Tertiary stats (column S)
$tertiaryRows = @(27,28,29,36,37)
foreach ($r in $tertiaryRows) {
    $ws.Range("S$r").Value = $true
}

# --- 5. New isMultiplier flags on rows that previously lacked one ---
$ws.Range("P32").Value = $true
$ws.Range("P34").Value = $true

# --- 6. Fill in the new "statMult" row (new row 38) ---
$ws.Range("A38").Value = "statMult"
$ws.Range("O38").Value = 1
$ws.Range("P38").Value = $true
$ws.Range("T38").Value = 1.4
$ws.Range("U38").Value = "Stat Multiplier: {0}"

# --- 7. Append the new "rarity" row at the end (row 42) ---
$ws.Range("A42").Value = "rarity"
$ws.Range("U42").Value = "Rarity: {0}"
